$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.688.64'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.024.66'
$ws.Range('E3').Value = '  +2.83%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.21'
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.90'
$ws.Range('E6').Value = '  +4.94%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.370'
$ws.Range('E11').Value = '  +5.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.545.39'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.36'
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000162'
$ws.Range('E15').Value = '  +3.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.704.12'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.028.03'
$ws.Range('E17').Value = '  +3.07%  '
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.13'
$ws.Range('E19').Value = '  +5.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.06'
$ws.Range('E20').Value = '  +4.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '333.95'
$ws.Range('E21').Value = '  +5.95%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.501'
$ws.Range('E23').Value = '  +4.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.84'
$ws.Range('E24').Value = '  +3.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.153.54'
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.167'
$ws.Range('E26').Value = '  +3.67%  '
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0925'
$ws.Range('E28').Value = '  +8.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.40'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('E31').Value = '  +3.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.43'
$ws.Range('E32').Value = '  +2.86%  '
$ws.Range('E33').Value = '  +3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '153.30'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.50'
$ws.Range('E35').Value = '  +1.63%  '
$ws.Range('B36').Value = 'EnergySwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '27.27'
$ws.Range('E36').Value = '  +15.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.82'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0664'
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.064.18'
$ws.Range('E40').Value = '  +3.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.51'
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +4.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.658'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.202.67'
$ws.Range('E45').Value = '  +3.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.35'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0246'
$ws.Range('E47').Value = '  +7.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.934'
$ws.Range('E48').Value = '  +3.33%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.85'
$ws.Range('E49').Value = '  +5.65%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.85'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0857'
$ws.Range('E51').Value = '  +1.52%  '
